$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (row 1) to snake_case field names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# 2) Title-case the Spanish connector words ("de", "del", "la", "el",
#    "los", "las", "y") inside state/municipality names, e.g.
#    "Pabellón de Arteaga" -> "Pabellón De Arteaga"
$ws.Range("B7").Value = 'Pabellón De Arteaga'
$ws.Range("B8").Value = 'Rincón De Romos'
$ws.Range("B28").Value = 'Amatenango De La Frontera'
$ws.Range("B32").Value = 'Benemérito De Las Américas'
$ws.Range("B34").Value = 'Comitán De Domínguez'
$ws.Range("B47").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B53").Value = 'Salto De Agua'
$ws.Range("B80").Value = 'Guadalupe Y Calvo'
$ws.Range("B83").Value = 'Hidalgo Del Parral'
$ws.Range("B94").Value = 'San Francisco De Borja'
$ws.Range("A98").Value = 'Ciudad De México'
$ws.Range("A115").Value = 'Coahuila De Zaragoza'
$ws.Range("B149").Value = 'Nombre De Dios'
$ws.Range("B153").Value = 'Pánuco De Coronado'
$ws.Range("B156").Value = 'San Juan Del Río'
$ws.Range("A163").Value = 'Estado De México'
$ws.Range("B163").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B165").Value = 'Almoloya De Juárez'
$ws.Range("B169").Value = 'Atizapán De Zaragoza'
$ws.Range("B177").Value = 'Ecatepec De Morelos'
$ws.Range("B180").Value = 'Ixtapan De La Sal'
$ws.Range("B184").Value = 'Naucalpan De Juárez'
$ws.Range("B189").Value = 'San Felipe Del Progreso'
$ws.Range("B197").Value = 'Tlalnepantla De Baz'
$ws.Range("B202").Value = 'Valle De Bravo'
$ws.Range("B203").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B206").Value = 'Villa Del Carbón'
$ws.Range("B211").Value = 'Apaseo El Alto'
$ws.Range("B216").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B224").Value = 'Purísima Del Rincón'
$ws.Range("B229").Value = 'San Diego De La Unión'
$ws.Range("B231").Value = 'San Francisco Del Rincón'
$ws.Range("B233").Value = 'San Luis De La Paz'
$ws.Range("B234").Value = 'San Miguel De Allende'
$ws.Range("B235").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B236").Value = 'Silao De La Victoria'
$ws.Range("B241").Value = 'Valle De Santiago'
$ws.Range("B244").Value = 'Acapulco De Juárez'
$ws.Range("B245").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B248").Value = 'Atenango Del Río'
$ws.Range("B250").Value = 'Atoyac De Álvarez'
$ws.Range("B251").Value = 'Ayutla De Los Libres'
$ws.Range("B254").Value = 'Chilapa De Álvarez'
$ws.Range("B255").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B258").Value = 'Coyuca De Benítez'
$ws.Range("B259").Value = 'Coyuca De Catalán'
$ws.Range("B263").Value = 'Cuetzala Del Progreso'
$ws.Range("B269").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B270").Value = 'Iguala De La Independencia'
$ws.Range("B271").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B286").Value = 'Taxco De Alarcón'
$ws.Range("B291").Value = 'Tlapa De Comonfort'
$ws.Range("B293").Value = 'Técpan De Galeana'
$ws.Range("B295").Value = 'Zihuatanejo De Azueta'
$ws.Range("B302").Value = 'Atotonilco De Tula'
$ws.Range("B304").Value = 'Cuautepec De Hinojosa'
$ws.Range("B307").Value = 'Huasca De Ocampo'
$ws.Range("B308").Value = 'Huejutla De Reyes'
$ws.Range("B312").Value = 'Mineral Del Chico'
$ws.Range("B313").Value = 'Mineral Del Monte'
$ws.Range("B314").Value = 'Mixquiahuala De Juárez'
$ws.Range("B315").Value = 'Pachuca De Soto'
$ws.Range("B317").Value = 'Progreso De Obregón'
$ws.Range("B319").Value = 'Tenango De Doria'
$ws.Range("B321").Value = 'Tepehuacán De Guerrero'
$ws.Range("B322").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B323").Value = 'Tezontepec De Aldama'
$ws.Range("B327").Value = 'Tula De Allende'
$ws.Range("B328").Value = 'Tulancingo De Bravo'
$ws.Range("B330").Value = 'Zacualtipán De Ángeles'
$ws.Range("B335").Value = 'Atemajac De Brizuela'
$ws.Range("B337").Value = 'Atotonilco El Alto'
$ws.Range("B338").Value = 'Autlán De Navarro'
$ws.Range("B352").Value = 'Huejuquilla El Alto'
$ws.Range("B358").Value = 'La Manzanilla De La Paz'
$ws.Range("B359").Value = 'Lagos De Moreno'
$ws.Range("B369").Value = 'San Diego De Alejandría'
$ws.Range("B370").Value = 'San Juan De Los Lagos'
$ws.Range("B372").Value = 'San Martín De Bolaños'
$ws.Range("B373").Value = 'San Miguel El Alto'
$ws.Range("B376").Value = 'Talpa De Allende'
$ws.Range("B377").Value = 'Tamazula De Gordiano'
$ws.Range("B379").Value = 'Tepatitlán De Morelos'
$ws.Range("B380").Value = 'Tizapán El Alto'
$ws.Range("B381").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B388").Value = 'Unión De San Antonio'
$ws.Range("B389").Value = 'Unión De Tula'
$ws.Range("B390").Value = 'Valle De Juárez'
$ws.Range("B394").Value = 'Yahualica De González Gallo'
$ws.Range("B398").Value = 'Zapotlán El Grande'
$ws.Range("A400").Value = 'Michoacán De Ocampo'
$ws.Range("B409").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B474").Value = 'Puente De Ixtla'
$ws.Range("B476").Value = 'Tlaltizapán De Zapata'
$ws.Range("B488").Value = 'Santa María Del Oro'
$ws.Range("B499").Value = 'San Nicolás De Los Garza'
$ws.Range("B504").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B508").Value = 'Fresnillo De Trujano'
$ws.Range("B509").Value = 'Guadalupe De Ramírez'
$ws.Range("B510").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B511").Value = 'Huajuapan De León'
$ws.Range("B515").Value = 'Mariscala De Juárez'
$ws.Range("B517").Value = 'Oaxaca De Juárez'
$ws.Range("B518").Value = 'Ocotlán De Morelos'
$ws.Range("B519").Value = 'Putla Villa De Guerrero'
$ws.Range("B522").Value = 'San Antonino El Alto'
$ws.Range("B536").Value = 'San Martín De Los Cansecos'
$ws.Range("B549").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B572").Value = 'Tataltepec De Valdés'
$ws.Range("B573").Value = 'Villa Sola De Vega'
$ws.Range("B575").Value = 'Zimatlán De Álvarez'
$ws.Range("B588").Value = 'Cuayuca De Andrade'
$ws.Range("B593").Value = 'Huehuetlán El Chico'
$ws.Range("B594").Value = 'Izúcar De Matamoros'
$ws.Range("B597").Value = 'Los Reyes De Juárez'
$ws.Range("B609").Value = 'Tepexi De Rodríguez'
$ws.Range("B614").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B624").Value = 'Amealco De Bonfil'
$ws.Range("B625").Value = 'Cadereyta De Montes'
$ws.Range("B628").Value = 'Jalpan De Serra'
$ws.Range("B631").Value = 'San Juan Del Río'
$ws.Range("B646").Value = 'Mexquitic De Carmona'
$ws.Range("B650").Value = 'Santa María Del Río'
$ws.Range("B654").Value = 'Villa De Arista'
$ws.Range("B655").Value = 'Villa De Ramos'
$ws.Range("B656").Value = 'Villa De La Paz'
$ws.Range("B709").Value = 'Soto La Marina'
$ws.Range("B716").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("A721").Value = 'Veracruz De Ignacio De La Llave'
$ws.Range("B723").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B726").Value = 'Castillo De Teayo'
$ws.Range("B740").Value = 'Hueyapan De Ocampo'
$ws.Range("B747").Value = 'Las Vigas De Ramírez'
$ws.Range("B748").Value = 'Martínez De La Torre'
$ws.Range("B756").Value = 'Ozuluama De Mascareñas'
$ws.Range("B758").Value = 'Paso De Ovejas'
$ws.Range("B762").Value = 'Poza Rica De Hidalgo'
$ws.Range("B769").Value = 'Soledad De Doblado'
$ws.Range("B793").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B803").Value = 'Jiménez Del Teul'
$ws.Range("B809").Value = 'Nochistlán De Mejía'
$ws.Range("B810").Value = 'Noria De Ángeles'
$ws.Range("B821").Value = 'Villa De Cos'

# 3) Normalize the grand-total label casing
$ws.Range("A825").Value = 'Total'

# 4) Remove the footer/metadata rows (827-831) describing sample size,
#    source and author notes that are no longer part of the clean table.
$ws.Range("A827:A831").EntireRow.Delete()

Write-Host ("Done. UsedRange=" + $ws.UsedRange.Address())
